# Generate Report for Handoff
# Replaces the previous handoff round's generated file names / timestamps
# with the new round's values, and appends two new rows (one per new
# dependency .png asset) to the Overview sheet and each language sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Constants shared across sheets
# ---------------------------------------------------------------------
$mdOld   = "f2a5d54c-cee0-42fd-9bb2-4ec3cd47035b.md"
$mdNew   = "6d1ed732-bde1-4be4-b012-04e94fa9abd7.md"
$png1    = "9a1eb7f0-9011-4d4d-92e4-6d945d17ccfd.png"
$png2    = "d923192f-9efd-4530-acea-2454b01703e1.png"

$ready   = "Ready for handoff"
$dtCommon = "2016-03-24 19:22:42"
$dtZh     = "2016-03-24 19:22:37"
$epoch    = "0001-01-01 00:00:00"
$include  = "Include"
$isDep    = "IsDependency"
$depFrom  = "e2e\6d1ed732-bde1-4be4-b012-04e94fa9abd7.md"

$zhXlfNew = "6d1ed732-bde1-4be4-b012-04e94fa9abd7.652c99d04754b2880d74136ce1728e0854086083.zh-cn.xlf"
$deXlfNew = "6d1ed732-bde1-4be4-b012-04e94fa9abd7.652c99d04754b2880d74136ce1728e0854086083.de-de.xlf"
$png1Target = "50005a7b1e167c947c1dccdd0fda1272d0e64329.png"
$png2Target = "d93a5a0dab01f06b6d500d92affc747048b0874d.png"

# Hyperlink target URL builders (keep the same repo / host patterns the
# workbook already used for the previous handoff round).
function Md-Url($fname) {
    return "https://github.com/OpenLocalizationTest/oltest/blob/245e3b38a197b69aa879b789652435fc455d0b01/e2e/" + $fname
}
function Zh-Url($fname) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e7e7cb0629a26b517a0a946c0f024cc396b6519b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $fname
}
function De-Url($fname) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b696994cf2319c7deee899780aad3af6c2e4102c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $fname
}

# Cornflower blue (FF6495ED) - same color as the workbook's existing
# "HyperLink" cell style - applied to newly-created hyperlink cells so
# they visually match the pre-existing hyperlink (A2) on every sheet.
$hyperlinkColor = 15570276

function Set-HyperlinkCell($ws, $addr, $text, $url) {
    $rng = $ws.Range($addr)
    $rng.Value = $text
    if ($rng.Hyperlinks.Count -gt 0) {
        $rng.Hyperlinks.Delete()
    }
    $ws.Hyperlinks.Add($rng, $url, "", "", $text) | Out-Null
    $rng.Font.Underline = 2
    $rng.Font.Color = $hyperlinkColor
}

function Set-DateCell($ws, $addr, $value) {
    $rng = $ws.Range($addr)
    $rng.Value = $value
    $rng.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Sheets.Item("Overview")

Set-HyperlinkCell $wsOverview "A2" $mdNew (Md-Url $mdNew)
$wsOverview.Range("B2").Value = $ready
$wsOverview.Range("C2").Value = $ready
Set-DateCell $wsOverview "D2" $dtCommon

Set-HyperlinkCell $wsOverview "A3" $png1 (Md-Url $png1)
$wsOverview.Range("B3").Value = $ready
$wsOverview.Range("C3").Value = $ready
Set-DateCell $wsOverview "D3" $dtCommon

Set-HyperlinkCell $wsOverview "A4" $png2 (Md-Url $png2)
$wsOverview.Range("B4").Value = $ready
$wsOverview.Range("C4").Value = $ready
Set-DateCell $wsOverview "D4" $dtCommon

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Sheets.Item("zh-cn")

Set-HyperlinkCell $wsZh "A2" $mdNew (Md-Url $mdNew)
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = $ready
Set-HyperlinkCell $wsZh "D2" $zhXlfNew (Zh-Url $zhXlfNew)
Set-DateCell $wsZh "E2" $dtZh
Set-DateCell $wsZh "H2" $epoch
$wsZh.Range("J2").Value = $include

$wsZh.Range("A3").Value = $png1
$wsZh.Range("B3").Value = ".png"
$wsZh.Range("C3").Value = $ready
$wsZh.Range("D3").Value = $png1Target
Set-DateCell $wsZh "E3" $dtZh
Set-DateCell $wsZh "H3" $epoch
$wsZh.Range("J3").Value = $isDep
$wsZh.Range("K3").Value = $depFrom
Set-HyperlinkCell $wsZh "A3" $png1 (Md-Url $png1)
Set-HyperlinkCell $wsZh "D3" $png1Target (Zh-Url $png1Target)

$wsZh.Range("A4").Value = $png2
$wsZh.Range("B4").Value = ".png"
$wsZh.Range("C4").Value = $ready
$wsZh.Range("D4").Value = $png2Target
Set-DateCell $wsZh "E4" $dtZh
Set-DateCell $wsZh "H4" $epoch
$wsZh.Range("J4").Value = $isDep
$wsZh.Range("K4").Value = $depFrom
Set-HyperlinkCell $wsZh "A4" $png2 (Md-Url $png2)
Set-HyperlinkCell $wsZh "D4" $png2Target (Zh-Url $png2Target)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Sheets.Item("de-de")

Set-HyperlinkCell $wsDe "A2" $mdNew (Md-Url $mdNew)
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = $ready
Set-HyperlinkCell $wsDe "D2" $deXlfNew (De-Url $deXlfNew)
Set-DateCell $wsDe "E2" $dtCommon
Set-DateCell $wsDe "H2" $epoch
$wsDe.Range("J2").Value = $include

$wsDe.Range("A3").Value = $png1
$wsDe.Range("B3").Value = ".png"
$wsDe.Range("C3").Value = $ready
$wsDe.Range("D3").Value = $png1Target
Set-DateCell $wsDe "E3" $dtCommon
Set-DateCell $wsDe "H3" $epoch
$wsDe.Range("J3").Value = $isDep
$wsDe.Range("K3").Value = $depFrom
Set-HyperlinkCell $wsDe "A3" $png1 (Md-Url $png1)
Set-HyperlinkCell $wsDe "D3" $png1Target (De-Url $png1Target)

$wsDe.Range("A4").Value = $png2
$wsDe.Range("B4").Value = ".png"
$wsDe.Range("C4").Value = $ready
$wsDe.Range("D4").Value = $png2Target
Set-DateCell $wsDe "E4" $dtCommon
Set-DateCell $wsDe "H4" $epoch
$wsDe.Range("J4").Value = $isDep
$wsDe.Range("K4").Value = $depFrom
Set-HyperlinkCell $wsDe "A4" $png2 (Md-Url $png2)
Set-HyperlinkCell $wsDe "D4" $png2Target (De-Url $png2Target)
